$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 7: Inscritos, Pagos, Inscrições homologadas updated (+1 each)
$ws.Range("E7").Value = 35
$ws.Range("F7").Value = 21
$ws.Range("H7").Value = 22

# Row 11: Inscritos updated (+1)
$ws.Range("E11").Value = 23

# Row 15: Inscritos, Pagos, Inscrições homologadas updated (+1 each)
$ws.Range("E15").Value = 115
$ws.Range("F15").Value = 53
$ws.Range("H15").Value = 64

# Row 18: Inscritos updated (+1)
$ws.Range("E18").Value = 102
